$wb = $excel.ActiveWorkbook

# --- Sheet "mean" ---
$wsMean = $wb.Worksheets.Item("mean")
$wsMean.Range("J3").Value = 3.97
$wsMean.Range("K3").Value = 0.088
$wsMean.Range("J4").Value = 3.952
$wsMean.Range("K4").Value = 0.106

# --- Sheet "stdev" ---
$wsStdev = $wb.Worksheets.Item("stdev")
$wsStdev.Range("J3").Value = 0.017
$wsStdev.Range("K3").Value = 0.017
$wsStdev.Range("J4").Value = 0.021
$wsStdev.Range("K4").Value = 0.021

# --- Sheet "summary" ---
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("L3").Value = 3.97
$wsSummary.Range("M3").Value = 0.017
$wsSummary.Range("T3").Value = 0.088
$wsSummary.Range("U3").Value = 0.017
$wsSummary.Range("L4").Value = 3.952
$wsSummary.Range("M4").Value = 0.021
$wsSummary.Range("T4").Value = 0.106
$wsSummary.Range("U4").Value = 0.021
